# Add a merged cell block (D1:E2) containing the text "merged cell" -
# mirrors the fixture update that adds a merged-cell case to the
# read_only worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Touch a format property on the anchor cell so it (and the other cells
# covered by the merge) are materialized with an explicit style in the
# saved worksheet instead of being dropped as all-default blanks.
$ws.Range("D1").Locked = $false

$ws.Range("D1").Value = "merged cell"

$ws.Range("D1:E2").Merge()

$ws.Range("D3").Select()
